$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from existing date/bool cells in row 3 down to row 5
$ws.Range("A3").Copy()
$ws.Range("A5").PasteSpecial(-4122)

$ws.Range("G3").Copy()
$ws.Range("G5").PasteSpecial(-4122)

# Now set the values for the new row
$ws.Range("A5").Value = 42635.643287037034
$ws.Range("B5").Value = $false
$ws.Range("C5").Value = 9931.5
$ws.Range("D5").Value = 10000
$ws.Range("E5").Value = 18.91
$ws.Range("F5").Value = 19.170000000000002
$ws.Range("G5").Value = $true
$ws.Range("H5").Value = 1.37
$ws.Range("I5").Value = $false
